# Insert a new column before column A to hold a row "ID" label,
# shifting the existing A:E data to B:F, and populate the new
# column with the ID header and per-row identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at position A; existing columns A:E shift to B:F.
$ws.Columns.Item(1).Insert()

# Copy the header formatting (bold/centered/bordered, style index 1)
# from the neighboring header cell (now B1) onto the new A1 header cell.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Set the new header text.
$ws.Range("A1").Value = "ID"

# Populate the ID values for each data row (rows 2-25).
$ids = @(
  "Hb 2",
  "Hb 3",
  "S 24",
  "S 28",
  "Hb 107",
  "Hb 66",
  "Hb 69",
  "Hb 95",
  "Hb 99",
  "Hb 92",
  "Hb 40",
  "Hb 41",
  "S 11",
  "Hb 57",
  "S 21",
  "S 22",
  "S 3",
  "S 4",
  "S 5",
  "Hb 74",
  "Hb 79",
  "Hb 32",
  "S 15",
  "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
  $ws.Cells.Item($i + 2, 1).Value = $ids[$i]
}
